# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 38 (date serial 44635,
# i.e. 2022-03-15), pushing the existing rows 38..284 down to 39..285.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 38; everything from the old row 38
# downward (through the old row 284) shifts down by one row.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly observation. The
# non-date fields mirror the row that used to sit at position 38 (same
# market / category / quality / prices / units / origin / classification);
# only the date (column D) is new.
$ws.Range("A38").Value = 3
$ws.Range("B38").Value = "Femacal de La Calera"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44635
$ws.Range("E38").Value = 5
$ws.Range("F38").Value = 100112039
$ws.Range("G38").Value = "Ciboulette"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 160
$ws.Range("K38").Value = 1500
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 1500
$ws.Range("N38").Value = "$/docena de atados"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 500
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = "Hortaliza"
